# Apply the "output generated at 456a3b4" gh-pages data refresh.
#
# The workbook has 4 sheets:
#   展览     (sheet1) - exhibitions
#   演出     (sheet2) - performances
#   本地生活 (sheet3) - local-life events
#   全部类型 (sheet4) - all types combined (cached aggregate of the above)
#
# For each, column F = 想去人数 (interested-count) and column G = 最低票价
# (lowest ticket price). This refresh bumps several interested-counts and
# flips a couple of now-sold-out listings' price cell from a number to the
# literal text "不可售" (not for sale).

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G3").Value = "不可售"
$ws.Range("F6").Value = 7525
$ws.Range("F8").Value = 7725
$ws.Range("F9").Value = 27
$ws.Range("F11").Value = 23
$ws.Range("F12").Value = 6355
$ws.Range("F13").Value = 3307
$ws.Range("F15").Value = 3668
$ws.Range("F16").Value = 31
$ws.Range("F17").Value = 27
$ws.Range("F18").Value = 28
$ws.Range("F19").Value = 43
$ws.Range("F21").Value = 448
$ws.Range("F23").Value = 296
$ws.Range("F25").Value = 3709
$ws.Range("F27").Value = 352
$ws.Range("F28").Value = 940
$ws.Range("F30").Value = 1361
$ws.Range("F31").Value = 68
$ws.Range("F32").Value = 32
$ws.Range("F33").Value = 2669
$ws.Range("F34").Value = 1647
$ws.Range("F35").Value = 22
$ws.Range("F38").Value = 3436
$ws.Range("F39").Value = 225
$ws.Range("F40").Value = 262
$ws.Range("F42").Value = 905
$ws.Range("F43").Value = 501
$ws.Range("F44").Value = 1333
$ws.Range("F47").Value = 607

# --- 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 68
$ws.Range("G3").Value = "不可售"

# --- 本地生活 (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 126

# --- 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 68
$ws.Range("G4").Value = "不可售"
$ws.Range("F6").Value = 126
$ws.Range("F11").Value = 7525
$ws.Range("F12").Value = 7725
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 23
$ws.Range("F15").Value = 6355
$ws.Range("F16").Value = 3307
$ws.Range("F18").Value = 3668
$ws.Range("F19").Value = 27
$ws.Range("F20").Value = 43
$ws.Range("F22").Value = 448
$ws.Range("F24").Value = 296
$ws.Range("F26").Value = 3709
$ws.Range("F29").Value = 352
$ws.Range("F30").Value = 940
$ws.Range("F32").Value = 1362
$ws.Range("F33").Value = 68
$ws.Range("F34").Value = 32
$ws.Range("F35").Value = 2669
$ws.Range("F36").Value = 1647
$ws.Range("F37").Value = 22
$ws.Range("F40").Value = 3436
$ws.Range("F41").Value = 225
$ws.Range("F42").Value = 262
$ws.Range("F44").Value = 905
$ws.Range("F45").Value = 501
$ws.Range("F46").Value = 1333
$ws.Range("F49").Value = 607
